$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new row at position 2 for the "PRODUTO" line. This
#    pushes the old PRODUTO row -> row3 (becomes "COR"), QTDE -> row4,
#    CODIGO -> row5. The B3/B4/B5 formatting (60pt/150pt/90pt wrap
#    styles) already land exactly where we need them after the shift,
#    so only the brand-new row 2 needs formatting copied in.
# ------------------------------------------------------------------
$ws.Rows.Item(2).Insert()

# New row-2 cells: copy formatting from row 3 (old PRODUTO row), which
# already carries the border/alignment/wrap/font we want to start from.
$ws.Range("A3:B3").Copy()
$ws.Range("A2:B2").PasteSpecial(-4122)  # xlPasteFormats

# A2 ("PRODUTO" label) gets its own, smaller bold font.
$ws.Range("A2").Font.Size = 24

# A3 ("COR", was "PRODUTO") label font shrinks 30 -> 28.
$ws.Range("A3").Font.Size = 28

# ------------------------------------------------------------------
# 2. Cell text content
# ------------------------------------------------------------------
$ws.Range("A2").Value = "PRODUTO"
$ws.Range("B2").Value = "Fita Borda Papel"
$ws.Range("A3").Value = "COR"
$ws.Range("B3").Value = "TITANIO"
$ws.Range("A4").Value = "QTDE"
$ws.Range("B4").Value = "'65320"
$ws.Range("A5").Value = "CODIGO"
$ws.Range("B5").Value = "'150001004"

# ------------------------------------------------------------------
# 3. Row heights
# ------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 32.25
$ws.Rows.Item(2).RowHeight = 120
$ws.Rows.Item(3).RowHeight = 138.75
$ws.Rows.Item(4).RowHeight = 225
$ws.Rows.Item(5).RowHeight = 236.25

# ------------------------------------------------------------------
# 4. Print area / dimension
# ------------------------------------------------------------------
$ws.PageSetup.PrintArea = "A1:B5"

"done"
